$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to Text format individually (per-cell) before assigning the value,
# since applying NumberFormat to a combined multi-area Range is unreliable here.
$textCells = @("D5","D7","D8","D9","D10","D11","D12","D13","D17","D18","D19","D20","D21","D22","D24","D25","D27","D28","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.141.81"
$ws.Range("E2").Value = "  +5.69%  "
$ws.Range("D3").Value = "1.921.73"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("D5").Value = "330.24"
$ws.Range("E5").Value = "  +4.73%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").Value = "0.5232"
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("D8").Value = "0.4094"
$ws.Range("E8").Value = "  +4.84%  "
$ws.Range("D9").Value = "0.08540"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "43.07"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").Value = "1.130"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").Value = "22.41"
$ws.Range("E12").Value = "  +9.48%  "
$ws.Range("D13").Value = "6.423"
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("D14").Value = "1.915.25"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "95.53"
$ws.Range("E17").Value = "  +4.71%  "
$ws.Range("D18").Value = "0.00001116"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "0.06690"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "18.47"
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "6.018"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").Value = "30.143.32"
$ws.Range("D24").Value = "11.40"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").Value = "2.214"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "2.134.72"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "160.25"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "21.11"
$ws.Range("D29").Value = "2.452"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").Value = "129.26"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("E31").Value = "  +3.57%  "
$ws.Range("D32").Value = "0.1065"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").Value = "6.056"
$ws.Range("E33").Value = "  +5.57%  "
$ws.Range("D34").Value = "3.638"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").Value = "0.02498"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "0.06617"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "0.2211"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").Value = "1.236"
$ws.Range("E38").Value = "  +4.59%  "
$ws.Range("D39").Value = "5.189"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "8.935"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "0.6555"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("D42").Value = "1.251"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").Value = "11.64"
$ws.Range("E43").Value = "  +4.84%  "
$ws.Range("D44").Value = "0.6170"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("D45").Value = "13.23"
$ws.Range("D46").Value = "3.761"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").Value = "2.082"
$ws.Range("E47").Value = "  +3.49%  "
$ws.Range("D48").Value = "1.252"
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("D49").Value = "124.76"
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").Value = "1.172"
$ws.Range("E50").Value = "  +11.36%  "
$ws.Range("D51").Value = "79.87"
$ws.Range("E51").Value = "  +4.22%  "
